$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-236). Update every occurrence of 45172 to 45175.
$ws.Range("C2:C236").Value = 45175
